$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 60 and row 61 contents (columns F:V) ---
# Read row 60 values (F..V => columns 6..22)
$f60 = $ws.Cells.Item(60, 6).Value2
$g60 = $ws.Cells.Item(60, 7).Value2
$h60 = $ws.Cells.Item(60, 8).Value2
$i60 = $ws.Cells.Item(60, 9).Value2
$j60 = $ws.Cells.Item(60, 10).Value2
$k60 = $ws.Cells.Item(60, 11).Value2
$l60 = $ws.Cells.Item(60, 12).Value2
$m60 = $ws.Cells.Item(60, 13).Value2
$n60 = $ws.Cells.Item(60, 14).Value2
$o60 = $ws.Cells.Item(60, 15).Value2
$p60 = $ws.Cells.Item(60, 16).Value2
$q60 = $ws.Cells.Item(60, 17).Value2
$r60 = $ws.Cells.Item(60, 18).Value2
$s60 = $ws.Cells.Item(60, 19).Value2
$t60 = $ws.Cells.Item(60, 20).Value2
$u60 = $ws.Cells.Item(60, 21).Value2
$v60 = $ws.Cells.Item(60, 22).Value2

# Read row 61 values (F..V => columns 6..22)
$f61 = $ws.Cells.Item(61, 6).Value2
$g61 = $ws.Cells.Item(61, 7).Value2
$h61 = $ws.Cells.Item(61, 8).Value2
$i61 = $ws.Cells.Item(61, 9).Value2
$j61 = $ws.Cells.Item(61, 10).Value2
$k61 = $ws.Cells.Item(61, 11).Value2
$l61 = $ws.Cells.Item(61, 12).Value2
$m61 = $ws.Cells.Item(61, 13).Value2
$n61 = $ws.Cells.Item(61, 14).Value2
$o61 = $ws.Cells.Item(61, 15).Value2
$p61 = $ws.Cells.Item(61, 16).Value2
$q61 = $ws.Cells.Item(61, 17).Value2
$r61 = $ws.Cells.Item(61, 18).Value2
$s61 = $ws.Cells.Item(61, 19).Value2
$t61 = $ws.Cells.Item(61, 20).Value2
$u61 = $ws.Cells.Item(61, 21).Value2
$v61 = $ws.Cells.Item(61, 22).Value2

# Write row 61's old values into row 60
$ws.Cells.Item(60, 6).Value2 = $f61
$ws.Cells.Item(60, 7).Value2 = $g61
$ws.Cells.Item(60, 8).Value2 = $h61
$ws.Cells.Item(60, 9).Value2 = $i61
$ws.Cells.Item(60, 10).Value2 = $j61
$ws.Cells.Item(60, 11).Value2 = $k61
$ws.Cells.Item(60, 12).Value2 = $l61
$ws.Cells.Item(60, 13).Value2 = $m61
$ws.Cells.Item(60, 14).Value2 = $n61
$ws.Cells.Item(60, 15).Value2 = $o61
$ws.Cells.Item(60, 16).Value2 = $p61
$ws.Cells.Item(60, 17).Value2 = $q61
$ws.Cells.Item(60, 18).Value2 = $r61
$ws.Cells.Item(60, 19).Value2 = $s61
$ws.Cells.Item(60, 20).Value2 = $t61
$ws.Cells.Item(60, 21).Value2 = $u61
$ws.Cells.Item(60, 22).Value2 = $v61

# Write row 60's old values into row 61
$ws.Cells.Item(61, 6).Value2 = $f60
$ws.Cells.Item(61, 7).Value2 = $g60
$ws.Cells.Item(61, 8).Value2 = $h60
$ws.Cells.Item(61, 9).Value2 = $i60
$ws.Cells.Item(61, 10).Value2 = $j60
$ws.Cells.Item(61, 11).Value2 = $k60
$ws.Cells.Item(61, 12).Value2 = $l60
$ws.Cells.Item(61, 13).Value2 = $m60
$ws.Cells.Item(61, 14).Value2 = $n60
$ws.Cells.Item(61, 15).Value2 = $o60
$ws.Cells.Item(61, 16).Value2 = $p60
$ws.Cells.Item(61, 17).Value2 = $q60
$ws.Cells.Item(61, 18).Value2 = $r60
$ws.Cells.Item(61, 19).Value2 = $s60
$ws.Cells.Item(61, 20).Value2 = $t60
$ws.Cells.Item(61, 21).Value2 = $u60
$ws.Cells.Item(61, 22).Value2 = $v60

# --- Append new row 82 ---
$ws.Cells.Item(82, 1).Value2 = 81
$ws.Cells.Item(82, 2).Value2 = "portugal"
$ws.Cells.Item(82, 3).Value2 = "liga-portugal-2"
$ws.Cells.Item(82, 4).Value2 = "2023-2024"
$ws.Cells.Item(82, 5).Value2 = 45236.79166666666
$ws.Cells.Item(82, 6).Value2 = "Pacos Ferreira"
$ws.Cells.Item(82, 7).Value2 = 1
$ws.Cells.Item(82, 8).Value2 = "Academico Viseu"
$ws.Cells.Item(82, 9).Value2 = 0
$ws.Cells.Item(82, 10).Value2 = 2.35
$ws.Cells.Item(82, 11).Value2 = "01/11/2023 16:12"
$ws.Cells.Item(82, 12).Value2 = 2.4
$ws.Cells.Item(82, 13).Value2 = "06/11/2023 18:58"
$ws.Cells.Item(82, 14).Value2 = 3.6
$ws.Cells.Item(82, 15).Value2 = "01/11/2023 16:12"
$ws.Cells.Item(82, 16).Value2 = 3.57
$ws.Cells.Item(82, 17).Value2 = "06/11/2023 18:58"
$ws.Cells.Item(82, 18).Value2 = 2.8
$ws.Cells.Item(82, 19).Value2 = "01/11/2023 16:12"
$ws.Cells.Item(82, 20).Value2 = 2.98
$ws.Cells.Item(82, 21).Value2 = "06/11/2023 18:58"
$ws.Cells.Item(82, 22).Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal-2/pacos-ferreira-academico-viseu/jZpnXCUu/"

# Apply same styles as row 81 for consistency (index column bold/border, date column format)
$ws.Cells.Item(81, 1).Copy()
$ws.Cells.Item(82, 1).PasteSpecial(-4122)

$ws.Cells.Item(81, 5).Copy()
$ws.Cells.Item(82, 5).PasteSpecial(-4122)
